$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the same date serial value (46061) for every
# data row (rows 2-120). The update bumps that date forward by one day
# (46061 -> 46062) for all of them.
for ($r = 2; $r -le 120; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46061) {
        $cell.Value2 = 46062
    }
}
